$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "43.915.50"

# Row 3
$ws.Range("D3").Value = "2.355.00"
$ws.Range("E3").Value = "  -0.56%  "

# Row 4
$ws.Range("E4").Value = "  +0.06%  "

# Row 5
$ws.Range("E5").Value = "  -3.18%  "

# Row 6
$ws.Range("D6").Value = "'240.96"
$ws.Range("D6").Style = "Normal"

# Row 7
$ws.Range("D7").Value = "'73.61"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.63%  "

# Row 8
$ws.Range("E8").Value = "  -0.05%  "

# Row 9
$ws.Range("E9").Value = "  +1.11%  "

# Row 10
$ws.Range("E10").Value = "  -2.44%  "

# Row 11
$ws.Range("D11").Value = "'59.29"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.23%  "

# Row 12
$ws.Range("D12").Value = "'33.78"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +6.15%  "

# Row 13
$ws.Range("D13").Value = "'7.35"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.85%  "

# Row 14
$ws.Range("E14").Value = "  -0.25%  "

# Row 15
$ws.Range("D15").Value = "2.705.51"
$ws.Range("E15").Value = "  -0.65%  "

# Row 16
$ws.Range("D16").Value = "'16.46"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.40%  "

# Row 17
$ws.Range("D17").Value = "'0.912"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.87%  "

# Row 18
$ws.Range("D18").Value = "2.353.56"
$ws.Range("E18").Value = "  -0.76%  "

# Row 19
$ws.Range("D19").Value = "43.820.34"
$ws.Range("E19").Value = "  -1.20%  "

# Row 21
$ws.Range("E21").Value = "  +0.44%  "

# Row 22
$ws.Range("D22").Value = "'77.86"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.52%  "

# Row 23
$ws.Range("D23").Value = "'257.46"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.48%  "

# Row 24
$ws.Range("E24").Value = "  +15.34%  "

# Row 25
$ws.Range("E25").Value = "  -0.02%  "

# Row 26
$ws.Range("E26").Value = "  -0.05%  "

# Row 27
$ws.Range("D27").Value = "'2.51"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.70%  "

# Row 28
$ws.Range("D28").Value = "'10.67"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.43%  "

# Row 29
$ws.Range("D29").Value = "'2.29"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.18%  "

# Row 30
$ws.Range("D30").Value = "'22.74"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.07%  "

# Row 31
$ws.Range("D31").Value = "'177.32"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.42%  "

# Row 32
$ws.Range("E32").Value = "  -0.14%  "

# Row 33
$ws.Range("E33").Value = "  +0.20%  "

# Row 34
$ws.Range("D34").Value = "'0.0757"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.71%  "

# Row 35
$ws.Range("E35").Value = "  -3.24%  "

# Row 36
$ws.Range("E36").Value = "  +2.10%  "

# Row 37
$ws.Range("E37").Value = "  -2.53%  "

# Row 38
$ws.Range("D38").Value = "'6.47"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.74%  "

# Row 39
$ws.Range("E39").Value = "  -3.97%  "

# Row 40
$ws.Range("E40").Value = "  +0.53%  "

# Row 41
$ws.Range("D41").Value = "'68.16"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +27.53%  "

# Row 42
$ws.Range("B42").Value = "Cronos"
$ws.Range("C42").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D42").Value = "'0.112"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +10.06%  "

# Row 43
$ws.Range("B43").Value = "FTXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D43").Value = "'5.12"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +15.05%  "

# Row 44
$ws.Range("D44").Value = "'9.32"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.29%  "

# Row 45
$ws.Range("E45").Value = "  +2.21%  "

# Row 46
$ws.Range("D46").Value = "'19.11"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.15%  "

# Row 47
$ws.Range("E47").Value = "  -0.06%  "

# Row 48
$ws.Range("E48").Value = "  +0.48%  "

# Row 50
$ws.Range("B50").Value = "Aave"
$ws.Range("C50").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D50").Value = "'99.83"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.46%  "

# Row 51
$ws.Range("B51").Value = "ARBITRUM"
$ws.Range("C51").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D51").Value = "'1.17"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.44%  "

